# Updated cryptos list (GitHub Actions style refresh of scraped
# coinranking.com data): refresh Price (column D) and Volume(1h)
# (column E) figures for the existing rows, and swap the
# RocketPoolETH / WEMIXToken rows (44 <-> 45) with their refreshed
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # The Price column stores plain text (e.g. "27.157.56", "214.25")
    # rather than numbers. Excel's COM layer auto-converts numeric-
    # looking strings (like "214.24") into real numbers when assigned
    # straight to .Value, so force a text number format first, write
    # the value, then drop the format back to the default "Normal"
    # style so no stray formatting is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '27.153.06'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').Value = '1.623.67'
$ws.Range('E4').Value = '  -0.05%  '
Set-TextCell $ws.Range('D5') '214.24'
$ws.Range('E5').Value = '  -1.42%  '
$ws.Range('E6').Value = '  +1.59%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -1.53%  '
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('E10').Value = '  +1.28%  '
Set-TextCell $ws.Range('D11') '0.0846'
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').Value = '1.628.31'
$ws.Range('E12').Value = '  -1.61%  '
$ws.Range('E13').Value = '  -0.39%  '
Set-TextCell $ws.Range('D14') '0.542'
$ws.Range('E14').Value = '  -0.56%  '
$ws.Range('D15').Value = '27.138.32'
$ws.Range('E15').Value = '  -0.35%  '
Set-TextCell $ws.Range('D16') '64.55'
$ws.Range('D17').Value = '0.0₃0745'
$ws.Range('E17').Value = '  +0.50%  '
Set-TextCell $ws.Range('D18') '215.74'
$ws.Range('E18').Value = '  -1.79%  '
$ws.Range('E19').Value = '  +0.00%  '
Set-TextCell $ws.Range('D20') '6.92'
$ws.Range('E20').Value = '  +0.43%  '
$ws.Range('E21').Value = '  -0.92%  '
Set-TextCell $ws.Range('D22') '2.40'
$ws.Range('E22').Value = '  -6.71%  '
$ws.Range('E23').Value = '  -2.00%  '
Set-TextCell $ws.Range('D24') '148.25'
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('E27').Value = '  -1.23%  '
$ws.Range('E28').Value = '  -1.17%  '
Set-TextCell $ws.Range('D29') '0.0507'
$ws.Range('E29').Value = '  -0.81%  '
$ws.Range('E30').Value = '  -1.31%  '
$ws.Range('E31').Value = '  -0.53%  '
Set-TextCell $ws.Range('D32') '3.00'
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('D33').Value = '1.343.91'
$ws.Range('E33').Value = '  +5.33%  '
$ws.Range('E34').Value = '  -0.56%  '
$ws.Range('E35').Value = '  -0.55%  '
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('E37').Value = '  +1.62%  '
$ws.Range('E38').Value = '  -0.18%  '
$ws.Range('E40').Value = '  -0.76%  '
Set-TextCell $ws.Range('D41') '65.62'
$ws.Range('E41').Value = '  +5.90%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('E43').Value = '  -1.43%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.760.82'
$ws.Range('E44').Value = '  -1.40%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell $ws.Range('D45') '0.928'
$ws.Range('E45').Value = '  +38.32%  '
Set-TextCell $ws.Range('D46') '90.03'
$ws.Range('E46').Value = '  -2.07%  '
$ws.Range('E47').Value = '  +0.82%  '
$ws.Range('E48').Value = '  -1.44%  '
Set-TextCell $ws.Range('D49') '0.100'
$ws.Range('E49').Value = '  +2.65%  '
$ws.Range('E50').Value = '  -0.59%  '
Set-TextCell $ws.Range('D51') '7.59'
$ws.Range('E51').Value = '  -1.47%  '
